$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: B1 becomes "number1" (new text), C1 becomes "number" (same formatting as old B1)
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B1").Value = "number1"
$ws.Range("C1").Value = "number"

# Fill in new data for C2 and C3 (C4:C6 left blank)
$ws.Range("C2").Value = 7878776426
$ws.Range("C3").Value = 7357240129

# Add new column C width, matching column A's bestFit width as closely as possible
$ws.Columns.Item(3).ColumnWidth = 10

# Update selection to D4
$ws.Range("D4").Select()

$wb.Save()
